$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells are formatted as text so numeric-looking strings
# (e.g. "56.350.21", "1.00") are preserved exactly as text, matching the
# source workbook which stores these as inline strings.
$cellValues = @{
    "D2" = "56.350.21"
    "E2" = "  +9.75%  "
    "D3" = "3.229.82"
    "E3" = "  +4.16%  "
    "E4" = "  +0.06%  "
    "D5" = "397.89"
    "E5" = "  +2.36%  "
    "D6" = "111.44"
    "E6" = "  +7.32%  "
    "E7" = "  +2.83%  "
    "D8" = "1.00"
    "E8" = "  -0.02%  "
    "E9" = "  +6.02%  "
    "D10" = "39.33"
    "E10" = "  +6.47%  "
    "D11" = "0.0915"
    "E11" = "  +6.78%  "
    "E12" = "  +2.04%  "
    "D13" = "3.739.98"
    "E13" = "  +4.31%  "
    "D14" = "8.13"
    "E14" = "  +5.01%  "
    "D15" = "19.10"
    "E15" = "  +3.04%  "
    "D16" = "3.230.59"
    "E16" = "  +4.28%  "
    "D17" = "1.05"
    "E17" = "  +4.81%  "
    "D18" = "10.93"
    "E18" = "  +1.84%  "
    "D19" = "56.239.06"
    "E19" = "  +9.33%  "
    "E20" = "  +3.22%  "
    "E21" = "  +6.74%  "
    "E22" = "  +4.52%  "
    "D23" = "298.83"
    "E23" = "  +12.29%  "
    "D24" = "75.87"
    "E24" = "  +8.44%  "
    "D25" = "3.22"
    "E25" = "  +1.54%  "
    "D26" = "8.16"
    "E26" = "  +1.91%  "
    "D27" = "28.16"
    "E27" = "  +2.51%  "
    "D28" = "7.46"
    "E28" = "  +2.80%  "
    "E29" = "  +4.50%  "
    "E30" = "  +0.35%  "
    "E31" = "  +4.18%  "
    "D32" = "11.16"
    "E32" = "  +7.17%  "
    "D33" = "0.0494"
    "E33" = "  +4.30%  "
    "D34" = "36.70"
    "E34" = "  +1.17%  "
    "D35" = "2.19"
    "E35" = "  +5.83%  "
    "D36" = "51.37"
    "E36" = "  +3.12%  "
    "D37" = "3.12"
    "E37" = "  +25.09%  "
    "D38" = "3.52"
    "E38" = "  +3.33%  "
    "D39" = "1.00"
    "E39" = "  +0.02%  "
    "D40" = "137.33"
    "E40" = "  +5.19%  "
    "D41" = "17.44"
    "E41" = "  +4.67%  "
    "E42" = "  +3.30%  "
    "E43" = "  +4.11%  "
    "B44" = "TheGraph"
    "C44" = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
    "D44" = "0.286"
    "E44" = "  -1.80%  "
    "B45" = "Stellar"
    "C45" = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
    "D45" = "0.119"
    "E45" = "  +2.93%  "
    "B46" = "ThetaToken"
    "C46" = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
    "D46" = "2.24"
    "E46" = "  +56.33%  "
    "B47" = "EnergySwap"
    "C47" = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
    "D47" = "22.22"
    "E47" = "  +0.20%  "
    "D48" = "2.46"
    "E48" = "  -2.03%  "
    "D49" = "2.09"
    "E49" = "  -0.13%  "
    "D50" = "2.130.82"
    "E50" = "  +2.63%  "
    "D51" = "0.0363"
    "E51" = "  +10.93%  "
}

foreach ($cellRef in $cellValues.Keys) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $cellValues[$cellRef]
}
